$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated roster: player name, position(s), team.
$data = @(
    @("Jalen Suggs", "PG,SG", "Orlando Magic"),
    @("Russell Westbrook", "PG,SG", "Denver Nuggets"),
    @("Jalen Green", "PG,SG", "Houston Rockets"),
    @("Chris Paul", "PG", "San Antonio Spurs"),
    @("Dillon Brooks", "SG,SF", "Houston Rockets"),
    @("Pascal Siakam", "SF,PF,C", "Indiana Pacers"),
    @("Jaylen Brown", "SG,SF", "Boston Celtics"),
    @("Khris Middleton", "SF", "Milwaukee Bucks"),
    @("Clint Capela", "C", "Atlanta Hawks"),
    @("Nikola Jokic", "C", "Denver Nuggets"),
    @("Jakob Poeltl", "C", "Toronto Raptors"),
    @("Dejounte Murray", "PG,SG", "New Orleans Pelicans"),
    @("Deni Avdija", "SF,PF", "Portland Trail Blazers"),
    @("Rudy Gobert", "C", "Minnesota Timberwolves"),
    @("Tyus Jones", "PG", "Phoenix Suns"),
    @("Paolo Banchero", "SF,PF", "Orlando Magic"),
    @("Chet Holmgren", "PF,C", "Oklahoma City Thunder"),
    @("Jerami Grant", "SF,PF", "Portland Trail Blazers")
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
    $ws.Cells.Item($row, 3).Value = $data[$i][2]
}
